$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.004
$ws1.Range("B2").Value = 0.969
$ws1.Range("C2").Value = 0.925

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = -0.009
$ws2.Range("B2").Value = 0.827
$ws2.Range("C2").Value = 0.919

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.026
$ws3.Range("B2").Value = 0.657
$ws3.Range("C2").Value = 2531.68
$ws3.Range("D2").Value = 0.02
